# Atualização do diário de bordo
# Applies the "Doc/DIÁRIO DE BORDO.xlsx" update:
#  - fixes row 31 (date/times)
#  - appends new logbook rows 32..43 with values, running totals and activity text
#  - widens column H and updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix existing row 31 (date moved from 15/03 to 16/03, times changed)
# ---------------------------------------------------------------------------
$ws.Cells.Item(31,1).Value2 = 43906
$ws.Cells.Item(31,2).Value2 = 0.375
$ws.Cells.Item(31,3).Value2 = 0.5

# ---------------------------------------------------------------------------
# Helper template: row 28 already has every column (A-H) populated with the
# correct cell styles (date / time / number / wrapped-text). Copying its
# formatting onto each freshly appended row lets Excel reuse the existing
# style table entries instead of inventing new ones.
# ---------------------------------------------------------------------------
$ws.Range("A28:H28").Copy()
$ws.Range("A32:H43").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Row 32
# ---------------------------------------------------------------------------
$ws.Cells.Item(32,1).Value2 = 43927
$ws.Cells.Item(32,2).Value2 = 0.54166666666666663
$ws.Cells.Item(32,3).Value2 = 0.66666666666666663
$ws.Cells.Item(32,4).Value2 = 0.79166666666666663
$ws.Cells.Item(32,5).Value2 = 0.875
$ws.Cells.Item(32,6).Value2 = 5
$ws.Cells.Item(32,7).Formula = "=G31+F32"
$ws.Cells.Item(32,8).Value2 = "Tentativa de responsividade"

# ---------------------------------------------------------------------------
# 3) Row 33
# ---------------------------------------------------------------------------
$ws.Cells.Item(33,1).Value2 = 43943
$ws.Cells.Item(33,2).Value2 = 0.375
$ws.Cells.Item(33,3).Value2 = 0.5
$ws.Cells.Item(33,4).Value2 = 0.75
$ws.Cells.Item(33,5).Value2 = 0.83333333333333337
$ws.Cells.Item(33,6).Value2 = 5
$ws.Cells.Item(33,7).Formula = "=G32+F33"
$ws.Cells.Item(33,8).Value2 = "Tentativa de responsividade"

# ---------------------------------------------------------------------------
# 4) Row 34
# ---------------------------------------------------------------------------
$ws.Cells.Item(34,1).Value2 = 43944
$ws.Cells.Item(34,2).Value2 = 0.375
$ws.Cells.Item(34,3).Value2 = 0.5
$ws.Cells.Item(34,4).Value2 = 0.70833333333333337
$ws.Cells.Item(34,5).Value2 = 0.79166666666666663
$ws.Cells.Item(34,6).Value2 = 5
$ws.Cells.Item(34,7).Formula = "=G33+F34"
$ws.Cells.Item(34,8).Value2 = "Tentativa de responsividade"

# ---------------------------------------------------------------------------
# 5) Row 35
# ---------------------------------------------------------------------------
$ws.Cells.Item(35,1).Value2 = 43945
$ws.Cells.Item(35,2).Value2 = 0.54166666666666663
$ws.Cells.Item(35,3).Value2 = 0.70833333333333337
$ws.Cells.Item(35,4).Value2 = 0.79166666666666663
$ws.Cells.Item(35,5).Value2 = 0.875
$ws.Cells.Item(35,6).Value2 = 6
$ws.Cells.Item(35,7).Formula = "=G34+F35"
$ws.Cells.Item(35,8).Value2 = "Finalização da responsividade"

# ---------------------------------------------------------------------------
# 6) Row 36
# ---------------------------------------------------------------------------
$ws.Cells.Item(36,1).Value2 = 44306
$ws.Cells.Item(36,2).Value2 = 0.33333333333333331
$ws.Cells.Item(36,3).Value2 = 0.45833333333333331
$ws.Cells.Item(36,4).Value2 = 0.75
$ws.Cells.Item(36,5).Value2 = 0.79166666666666663
$ws.Cells.Item(36,6).Value2 = 4
$ws.Cells.Item(36,7).Formula = "=G35+F36"
$ws.Cells.Item(36,8).Value2 = "Instalação e configurações - retorno"

# ---------------------------------------------------------------------------
# 7) Row 37 (no D/E times, G is a hard-coded value, not a formula)
# ---------------------------------------------------------------------------
$ws.Cells.Item(37,1).Value2 = 44316
$ws.Cells.Item(37,2).Value2 = 0.54861111111111105
$ws.Cells.Item(37,3).Value2 = 0.59027777777777779
$ws.Cells.Item(37,4).ClearContents()
$ws.Cells.Item(37,5).ClearContents()
$ws.Cells.Item(37,6).Value2 = 1
$ws.Cells.Item(37,7).Value2 = 152
$ws.Cells.Item(37,8).Value2 = "Reunião com o professor Leonardo"

# ---------------------------------------------------------------------------
# 8) Row 38
# ---------------------------------------------------------------------------
$ws.Cells.Item(38,1).Value2 = 44319
$ws.Cells.Item(38,2).Value2 = 0.33333333333333331
$ws.Cells.Item(38,3).Value2 = 0.41666666666666669
$ws.Cells.Item(38,4).Value2 = 0.70833333333333337
$ws.Cells.Item(38,5).Value2 = 0.83333333333333337
$ws.Cells.Item(38,6).Value2 = 5
$ws.Cells.Item(38,7).Formula = "=G37+F38"
$ws.Cells.Item(38,8).Value2 = "Desenvolvimento da página para o sistema"

# ---------------------------------------------------------------------------
# 9) Row 39
# ---------------------------------------------------------------------------
$ws.Cells.Item(39,1).Value2 = 44320
$ws.Cells.Item(39,2).Value2 = 0.33333333333333331
$ws.Cells.Item(39,3).Value2 = 0.45833333333333331
$ws.Cells.Item(39,4).Value2 = 0.75
$ws.Cells.Item(39,5).Value2 = 0.83333333333333337
$ws.Cells.Item(39,6).Value2 = 5
$ws.Cells.Item(39,7).Formula = "=G38+F39"
$ws.Cells.Item(39,8).Value2 = "Desenvolvimento da página para o sistema"

# ---------------------------------------------------------------------------
# 10) Row 40
# ---------------------------------------------------------------------------
$ws.Cells.Item(40,1).Value2 = 44321
$ws.Cells.Item(40,2).Value2 = 0.54166666666666663
$ws.Cells.Item(40,3).Value2 = 0.70833333333333337
$ws.Cells.Item(40,4).Value2 = 0.79166666666666663
$ws.Cells.Item(40,5).Value2 = 0.875
$ws.Cells.Item(40,6).Value2 = 6
$ws.Cells.Item(40,7).Formula = "=G39+F40"
$ws.Cells.Item(40,8).Value2 = "Desenvolvimento da página para o sistema"

# ---------------------------------------------------------------------------
# 11) Row 41 (no D/E times)
# ---------------------------------------------------------------------------
$ws.Cells.Item(41,1).Value2 = 44322
$ws.Cells.Item(41,2).Value2 = 0.54166666666666663
$ws.Cells.Item(41,3).Value2 = 0.625
$ws.Cells.Item(41,4).ClearContents()
$ws.Cells.Item(41,5).ClearContents()
$ws.Cells.Item(41,6).Value2 = 2
$ws.Cells.Item(41,7).Formula = "=G40+F41"
$ws.Cells.Item(41,8).Value2 = "Ajustes finais e conclusão da página para o sistema"

# ---------------------------------------------------------------------------
# 12) Row 42
# ---------------------------------------------------------------------------
$ws.Cells.Item(42,1).Value2 = 44323
$ws.Cells.Item(42,2).Value2 = 0.41666666666666669
$ws.Cells.Item(42,3).Value2 = 0.5
$ws.Cells.Item(42,4).Value2 = 0.54861111111111105
$ws.Cells.Item(42,5).Value2 = 0.59027777777777779
$ws.Cells.Item(42,6).Value2 = 3
$ws.Cells.Item(42,7).Formula = "=G41+F42"
$ws.Cells.Item(42,8).Value2 = "Tentativa de correção do erro na criação de usuário e reunião com o professor Leonardo"

# ---------------------------------------------------------------------------
# 13) Row 43
# ---------------------------------------------------------------------------
$ws.Cells.Item(43,1).Value2 = 44328
$ws.Cells.Item(43,2).Value2 = 0.33333333333333331
$ws.Cells.Item(43,3).Value2 = 0.45833333333333331
$ws.Cells.Item(43,4).Value2 = 0.66666666666666663
$ws.Cells.Item(43,5).Value2 = 0.75
$ws.Cells.Item(43,6).Value2 = 5
$ws.Cells.Item(43,7).Formula = "=G42+F43"
$ws.Cells.Item(43,8).Value2 = "Correção de erro na criação de usuário"

# ---------------------------------------------------------------------------
# Recalculate so every running-total formula shows its final value.
# ---------------------------------------------------------------------------
$wb.Application.Calculate()

# ---------------------------------------------------------------------------
# Cosmetic touch-ups: widen column H, scroll/select to match the saved view.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 77.2

$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("G40").Select()
